# Updates the cryptos price/volume table to the latest scraped snapshot.
# Numeric-looking "Price" cells are forced to Text (NumberFormat "@") before
# assignment so Excel doesn't silently convert them to floating point
# numbers (which would lose the original text formatting, e.g. "42.413.08"
# or introduce binary float rounding, e.g. 306.39 -> 306.389999...).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.413.08'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '2.275.58'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.39'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.13'
$ws.Range('E6').Value = '  +3.39%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.46'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.92'
$ws.Range('E13').Value = '  +3.50%  '
$ws.Range('D14').Value = '2.627.77'
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.82'
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('D16').Value = '2.278.83'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.797'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '42.265.34'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.56'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.28'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.36'
$ws.Range('E23').Value = '  -1.75%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.57'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.24'
$ws.Range('E28').Value = '  +5.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.56'
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.74'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.19'
$ws.Range('E34').Value = '  +3.41%  '
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.64'
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('E38').Value = '  -3.33%  '
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.10'
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('E42').Value = '  +2.31%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.952.63'
$ws.Range('E43').Value = '  -3.13%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.97'
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.96'
$ws.Range('E46').Value = '  -2.54%  '
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.88'
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '92.67'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.13'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('E51').Value = '  -1.65%  '
